# correcion en el powerpoint
# - Refresh the cached "today" date shown by the date placeholders (master,
#   every slide layout and the notes master) from 28/2/2022 to 1/3/2022.
# - On slide 5, widen the first flow-chart box and tweak its first line of
#   text ("Las reservas..." -> "De las reservas...").

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText {
    param(
        $Shapes,
        [string]$NewText
    )

    for ($i = 1; $i -le $Shapes.Count; $i++) {
        $shp = $Shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }

        if ($isDatePlaceholder -and $shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $shp.TextFrame.TextRange.Text = $NewText
        }
    }
}

$newDate = "1/3/2022"

# Slide master date placeholder.
$master = $p.SlideMaster
Set-DatePlaceholderText -Shapes $master.Shapes -NewText $newDate

# Every slide layout's date placeholder.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Set-DatePlaceholderText -Shapes $layout.Shapes -NewText $newDate
}

# Notes master date placeholder.
$notesMaster = $p.NotesMaster
Set-DatePlaceholderText -Shapes $notesMaster.Shapes -NewText $newDate

# Slide 5: widen the first "flowchart alternate process" box and update its
# first line of text.
$slide5 = $p.Slides.Item(5)
$shape = $slide5.Shapes.Item(1)
$shape.Width = 483.283035

$textRange = $shape.TextFrame.TextRange
$firstLine = $textRange.Characters(1, 42)
$firstLine.Text = "De las reservas se deberán poder visualizar: "
